# blankRowInserter.py was streamlined from two separate loops into one.
# Net effect on the sheet:
#   - 3 additional blank rows are inserted above the fruit data block
#     (it used to start at row 5, now it starts at row 8), and
#   - the sheet now ends with 5 trailing blank rows instead of 2
#     (rows 13-17 instead of rows 10-11).
# The already-populated header rows (1-2) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right above the fruit data (currently rows 5-9),
# pushing that block down to rows 8-12. EntireRow.Insert() shifts
# everything from row 3 downward and copies formatting onto the newly
# created rows, so...
$ws.Rows.Item(3).Resize(3).EntireRow.Insert()

# ...immediately wipe the 3 new rows clean (content + formatting) so they
# stay truly blank, just like the untouched rows 3-4 in the original file.
$ws.Range("A3:C5").Clear()

# Make sure the sheet now ends with 5 blank rows (13-17) rather than the
# original 2 (10-11, which got consumed by the insert above). Touching a
# single cell per row with a no-op style keeps each row present/blank
# while extending the sheet's dimension down to row 17.
for ($r = 13; $r -le 17; $r++) {
  $ws.Cells.Item($r, 1).Style = "Normal"
}
